# Applies the betexplorer scrape "re-run" update:
#  - 14 pairs of adjacent rows had their match data (cols F:V) swapped
#    (the A:E "meta" columns - Indice/pais/torneio/temporada/data_partida -
#    stayed put while the actual match rows that got scraped in a different
#    order were exchanged).
#  - 2 brand-new match rows were appended at the bottom (231/232).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchData {
    param($RowA, $RowB)

    $rangeA = $ws.Range("F$RowA`:V$RowA")
    $rangeB = $ws.Range("F$RowB`:V$RowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

$rowPairs = @(
    @(43, 44),
    @(59, 60),
    @(72, 73),
    @(82, 83),
    @(85, 86),
    @(130, 131),
    @(141, 142),
    @(152, 153),
    @(158, 159),
    @(163, 164),
    @(169, 170),
    @(173, 174),
    @(175, 176),
    @(184, 185)
)

foreach ($pair in $rowPairs) {
    Swap-MatchData $pair[0] $pair[1]
}

# Append the two new rows (231 and 232), copying formatting from the
# last existing data row (230) so styles (bold index column, date format
# column, etc.) stay consistent with the rest of the sheet.
$ws.Range("A230:V230").Copy()
$ws.Range("A231:V232").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @{
        Row = 231
        Indice = 230
        Home = "Amorebieta"; HomeGoals = 1; Away = "Alcorcon"; AwayGoals = 2
        HomeOpenOdds = 2.44; HomeOpenDt = "17/12/2023 18:43"
        HomeCloseOdds = 2.19; HomeCloseDt = "21/12/2023 21:27"
        DrawOpenOdds = 3.08; DrawOpenDt = "17/12/2023 18:43"
        DrawCloseOdds = 3.08; DrawCloseDt = "21/12/2023 21:27"
        AwayOpenOdds = 3.29; AwayOpenDt = "17/12/2023 18:43"
        AwayCloseOdds = 4.06; AwayCloseDt = "21/12/2023 21:27"
        Url = "https://www.betexplorer.com/football/spain/laliga2/amorebieta-alcorcon/fcdcqyIt/"
    },
    @{
        Row = 232
        Indice = 231
        Home = "Valladolid"; HomeGoals = 0; Away = "Ferrol"; AwayGoals = 1
        HomeOpenOdds = 1.92; HomeOpenDt = "18/12/2023 20:42"
        HomeCloseOdds = 1.88; HomeCloseDt = "21/12/2023 21:29"
        DrawOpenOdds = 3.33; DrawOpenDt = "18/12/2023 20:42"
        DrawCloseOdds = 3.51; DrawCloseDt = "21/12/2023 21:28"
        AwayOpenOdds = 4.54; AwayOpenDt = "18/12/2023 20:42"
        AwayCloseOdds = 4.71; AwayCloseDt = "21/12/2023 21:29"
        Url = "https://www.betexplorer.com/football/spain/laliga2/valladolid-ferrol/pnkAtZ1a/"
    }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value2 = $nr.Indice
    $ws.Cells.Item($r, 2).Value2 = "spain"
    $ws.Cells.Item($r, 3).Value2 = "laliga2"
    $ws.Cells.Item($r, 4).Value2 = "2023-2024"
    $ws.Cells.Item($r, 5).Value2 = 45281.89583333334
    $ws.Cells.Item($r, 6).Value2 = $nr.Home
    $ws.Cells.Item($r, 7).Value2 = $nr.HomeGoals
    $ws.Cells.Item($r, 8).Value2 = $nr.Away
    $ws.Cells.Item($r, 9).Value2 = $nr.AwayGoals
    $ws.Cells.Item($r, 10).Value2 = $nr.HomeOpenOdds
    $ws.Cells.Item($r, 11).Value2 = $nr.HomeOpenDt
    $ws.Cells.Item($r, 12).Value2 = $nr.HomeCloseOdds
    $ws.Cells.Item($r, 13).Value2 = $nr.HomeCloseDt
    $ws.Cells.Item($r, 14).Value2 = $nr.DrawOpenOdds
    $ws.Cells.Item($r, 15).Value2 = $nr.DrawOpenDt
    $ws.Cells.Item($r, 16).Value2 = $nr.DrawCloseOdds
    $ws.Cells.Item($r, 17).Value2 = $nr.DrawCloseDt
    $ws.Cells.Item($r, 18).Value2 = $nr.AwayOpenOdds
    $ws.Cells.Item($r, 19).Value2 = $nr.AwayOpenDt
    $ws.Cells.Item($r, 20).Value2 = $nr.AwayCloseOdds
    $ws.Cells.Item($r, 21).Value2 = $nr.AwayCloseDt
    $ws.Cells.Item($r, 22).Value2 = $nr.Url
}
